$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Enter new shared strings in the same order the author typed them, so that
# the sharedStrings table indices line up with the target workbook.
$ws.Range("A64").Value = "watch"
$ws.Range("B64").Value = "Watch a command repeatly"

$ws.Range("A65").Value = "tree"
$ws.Range("B65").Value = "List a directory with selected depth"
$ws.Range("C65").Value = "tree -L 2   //list the directory with depth 2"

$ws.Range("C64").Value = "watch -n 2 tree    // repeatly update tree command result to stdout by 2 second interval"

# Update the view so the newly added row is visible / selected, matching author's final state
$ws.Range("A65").Select()
$excel.ActiveWindow.ScrollRow = 60
$excel.ActiveWindow.ScrollColumn = 1
